# "combine all solutions in one folder" - append weeks 11 & 12 of tracking data,
# add a short-term planning table, clean up the old REVIEW-plan notes, and
# update the running totals/percentage formula accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Running-total row: add a manual "extra" bucket (F5) with its own note,
#    and widen the grand-total formula in F3 to include it.
# ---------------------------------------------------------------------------
$ws.Range("F5").Font.Size = 20
$ws.Range("F5").Interior.Color = 65535
$ws.Range("F5").Value = 26
$ws.Range("G5").Font.Size = 20
$ws.Range("G5").Value = "<- sql题 暂时忽略"

$ws.Range("F3").Formula = "=SUM(C3:E3,F5,F4)"

# ---------------------------------------------------------------------------
# 2. Tidy up the old "REVIEW plan" notes in G76:G78 - the plan was completed,
#    so the notes are cleared (formatting reset to the plain/no-fill look).
# ---------------------------------------------------------------------------
$ws.Range("A11").Copy()
$ws.Range("G76").PasteSpecial(-4122)
$ws.Range("G77").PasteSpecial(-4122)
$ws.Range("G78").PasteSpecial(-4122)
$ws.Range("G76").ClearContents()
$ws.Range("G77").ClearContents()
$ws.Range("G78").ClearContents()

# ---------------------------------------------------------------------------
# 3. Week 11 (row 84-90) - highlighted red, was the "catch up" week.
# ---------------------------------------------------------------------------
$ws.Range("A84:E90").Font.Size = 20
$ws.Range("A84:E90").Interior.Color = 255

$ws.Range("A84").Value = 11
$ws.Range("B84").Value = "Mon"
$ws.Range("D84").Value = 6

$ws.Range("A85").Value = 43262
$ws.Range("A85").NumberFormat = "m/d/yy"
$ws.Range("B85").Value = "Tue"
$ws.Range("D85").Value = 5

$ws.Range("B86").Value = "Wed"
$ws.Range("C86").Value = 2
$ws.Range("D86").Value = 8

$ws.Range("B87").Value = "Thu"
$ws.Range("C87").Value = 1
$ws.Range("D87").Value = 9

$ws.Range("B88").Value = "Fri"
$ws.Range("C88").Value = 1
$ws.Range("D88").Value = 4

$ws.Range("B89").Value = "Sat"
$ws.Range("C89").Value = 5
$ws.Range("D89").Value = 5

$ws.Range("B90").Value = "Sun"
$ws.Range("C90").Value = 3
$ws.Range("D90").Value = 2
$ws.Range("E90").Value = 5

$ws.Range("F90").Font.Size = 20
$ws.Range("F90").Interior.Color = 65535
$ws.Range("F90").Formula = "=SUM(C84:E90)"

# ---------------------------------------------------------------------------
# 4. Plan table (H84:K94) sitting next to weeks 11/12 - highlighted orange.
# ---------------------------------------------------------------------------
$ws.Range("H84:K94").Font.Size = 20
$ws.Range("H84:K94").Interior.Color = 49407

$ws.Range("I84").Value = "Plan for the next 3 weeks"
$ws.Range("K84").Value = "完成度"

$ws.Range("H85").Value = "week 1"
$ws.Range("I85").Value = 40
$ws.Range("K85").Value = "完成了35"

$ws.Range("H86").Value = "week 2"
$ws.Range("I86").Value = 75

$ws.Range("H87").Value = "week 3"
$ws.Range("I87").Value = 75
$ws.Range("J87").Value = "完成前400"

$ws.Range("I89").Value = "Plan for the rest 4 weeks"
$ws.Range("J89").Value = "average daily"

$ws.Range("H90").Value = "week 1"
$ws.Range("I90").Value = "复习1-150"
$ws.Range("J90").Value = 21

$ws.Range("H91").Value = "week 2"
$ws.Range("I91").Value = "复习151-300"
$ws.Range("J91").Value = 21

$ws.Range("H92").Value = "week 3"
$ws.Range("I92").Value = "复习300-400, 三刷1-100"
$ws.Range("J92").Value = 28

$ws.Range("H93").Value = "week 4"
$ws.Range("I93").Value = "三刷101-300"
$ws.Range("J93").Value = 28

$ws.Range("H94").Value = "在家10天"
$ws.Range("I94").Value = "三刷301-400"
$ws.Range("J94").Value = 10

# ---------------------------------------------------------------------------
# 5. Week 12 (row 92-98) - back to the normal (no fill) look.
# ---------------------------------------------------------------------------
$ws.Range("A92").Value = 12
$ws.Range("B92").Value = "Mon"
$ws.Range("C92").Value = 5
$ws.Range("E92").Value = 5

$ws.Range("A93").Value = 43269
$ws.Range("A93").NumberFormat = "m/d/yy"
$ws.Range("B93").Value = "Tue"
$ws.Range("C93").Value = 6
$ws.Range("D93").Value = 3
$ws.Range("E93").Value = 1

$ws.Range("B94").Value = "Wed"
$ws.Range("C94").Value = 3
$ws.Range("D94").Value = 4
$ws.Range("E94").Value = 3

$ws.Range("B95").Value = "Thu"
$ws.Range("D95").Value = 5

$ws.Range("B96").Value = "Fri"
$ws.Range("D96").Value = 8
$ws.Range("E96").Value = 2

$ws.Range("B97").Value = "Sat"
$ws.Range("D97").Value = 10

$ws.Range("B98").Value = "Sun"
$ws.Range("D98").Value = 2
$ws.Range("E98").Value = 1

$ws.Range("F98").Font.Size = 20
$ws.Range("F98").Interior.Color = 65535
$ws.Range("F98").Formula = "=SUM(C92:E98)"
$ws.Range("G98").Font.Size = 20
$ws.Range("G98").Value = "每天必须要做5道hard"
$ws.Range("G99").Font.Size = 20
$ws.Range("G99").Value = "剩下5道medium"

# ---------------------------------------------------------------------------
# 6. Week 12 continued (row 100-106) - a second tracking block for the
#    same week, still in progress.
# ---------------------------------------------------------------------------
$ws.Range("A100").Value = 12
$ws.Range("B100").Value = "Mon"
$ws.Range("D100").Value = 3

$ws.Range("A101").Value = 43276
$ws.Range("A101").NumberFormat = "m/d/yy"
$ws.Range("B101").Value = "Tue"

$ws.Range("B102").Value = "Wed"
$ws.Range("B103").Value = "Thu"
$ws.Range("B104").Value = "Fri"
$ws.Range("B105").Value = "Sat"
$ws.Range("B106").Value = "Sun"

$ws.Range("F106").Font.Size = 20
$ws.Range("F106").Interior.Color = 65535
$ws.Range("F106").Formula = "=SUM(C100:E106)"

# ---------------------------------------------------------------------------
# 7. View/window cosmetics to match where the author ended up.
# ---------------------------------------------------------------------------
$ws.Range("D100").Select()
$excel.ActiveWindow.Zoom = 136
$ws.Application.ActiveWindow.ScrollRow = 87
